$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row 2 (BMI): algorithm label "BMI0" -> "bmi0"
$ws.Range("F2").Value = "bmi0"

# Row 3 (ENERGY): algorithm / comment labels "gj" -> "GJ" (uppercase)
$ws.Range("F3").Value = "GJ"
$ws.Range("H3").Value = "GJ*4.184"
$ws.Range("I3").Value = "GJ [kJ] calculated in kcal"

# Row 12 (WAIST): "waist" -> "waist0"
$ws.Range("F12").Value = "waist0"

# Row 13 (HIP): "hip" -> "hip0"
$ws.Range("F13").Value = "hip0"

# Rows 20-27 (CARB, PROT, FAT, ALC, FIBER, SFA, MUFA, PUFA): lowercase -> uppercase
$ws.Range("F20").Value = "ZK"
$ws.Range("F21").Value = "ZE"
$ws.Range("F22").Value = "ZF"
$ws.Range("F23").Value = "ZA"
$ws.Range("F24").Value = "ZB"
$ws.Range("F25").Value = "FS"
$ws.Range("F26").Value = "FU"
$ws.Range("F27").Value = "FP"

# Row 28 (TOT_SUGARS): "kd;km" -> "KD;KM", "kd + km" -> "KD+KM"
$ws.Range("F28").Value = "KD;KM"
$ws.Range("H28").Value = "KD+KM"

# Rows 31-32 (GLUC, FRUC): lowercase -> uppercase
$ws.Range("F31").Value = "KMT"
$ws.Range("F32").Value = "KMF"

# Row 35 (SODIUM): "mna" -> "MNA"
$ws.Range("F35").Value = "MNA"

# Row 36 (SOD_POT_RATIO): "mna;mk" -> "MNA;MK", "mna/mk" -> "MNA/MK"
$ws.Range("F36").Value = "MNA;MK"
$ws.Range("H36").Value = "MNA/MK"

# Update the sheet view: zoom and selected cell
$excel.ActiveWindow.Zoom = 70
$ws.Range("E5").Select()
